$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C and D hold numeric-looking identifiers (e.g. "44331122") that must
# stay stored as text, matching the rest of the sheet's ID columns.
$ws.Range("C41:D41").NumberFormat = "@"

$ws.Range("A41").Value = "EZEZ"
$ws.Range("B41").Value = "DD"
$ws.Range("C41").Value = "44331122"
$ws.Range("D41").Value = "11223344"
$ws.Range("E41").Value = "fdfdd"
$ws.Range("F41").Value = "الماشية"
$ws.Range("G41").Value = "الأبقار"
